$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "38.228.99"
$ws.Range("E2").Value = "  +3.09%  "
$ws.Range("D3").Value = "2.062.79"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'230.47"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("D7").Value = "'59.71"
$ws.Range("E7").Value = "  +8.63%  "
$ws.Range("D9").Value = "'0.388"
$ws.Range("E9").Value = "  +3.41%  "
$ws.Range("D10").Value = "'0.0815"
$ws.Range("E10").Value = "  +4.48%  "
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D12").Value = "'14.77"
$ws.Range("E12").Value = "  +5.12%  "
$ws.Range("D13").Value = "2.362.65"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "'21.32"
$ws.Range("E14").Value = "  +7.72%  "
$ws.Range("D15").Value = "'0.757"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D17").Value = "2.074.73"
$ws.Range("E17").Value = "  +3.47%  "
$ws.Range("D18").Value = "38.093.25"
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("D19").Value = "'6.29"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").Value = "'70.00"
$ws.Range("E20").Value = "  +2.49%  "
$ws.Range("E21").Value = "  +3.13%  "
$ws.Range("D22").Value = "'225.01"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("E25").Value = "  +4.22%  "
$ws.Range("D26").Value = "'9.34"
$ws.Range("E26").Value = "  +4.15%  "
$ws.Range("D27").Value = "'166.30"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("E28").Value = "  +7.50%  "
$ws.Range("D29").Value = "'19.08"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("E30").Value = "  +3.23%  "
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").Value = "'4.59"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("D34").Value = "'2.06"
$ws.Range("E34").Value = "  +11.09%  "
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").Value = "'2.34"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "'6.14"
$ws.Range("E37").Value = "  +15.05%  "
$ws.Range("D38").Value = "'3.29"
$ws.Range("E38").Value = "  +5.60%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "1.533.86"
$ws.Range("E40").Value = "  +5.47%  "
$ws.Range("D41").Value = "'98.42"
$ws.Range("E41").Value = "  +3.79%  "
$ws.Range("E42").Value = "  +2.55%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'16.92"
$ws.Range("E43").Value = "  +6.09%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.88"
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").Value = "'4.10"
$ws.Range("E47").Value = "  +13.85%  "
$ws.Range("D48").Value = "'1.03"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").Value = "'2.98"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("D50").Value = "'7.12"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "2.250.10"
